# Insert two new data rows above the existing row 703 ("Femacal de La
# Calera" / Cebollín price sheet). Excel's row-insert shifts all the
# rows that used to be 703-727 down to 705-729, which matches the
# dimension growing from A1:R727 to A1:R729 in the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(703).EntireRow.Insert()
$ws.Rows.Item(703).EntireRow.Insert()

# New row 703: Cebollín, Primera, market date 2023-05-29
$ws.Cells.Item(703, 1).Value = 3
$ws.Cells.Item(703, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(703, 3).Value = 'Coquimbo'
$ws.Cells.Item(703, 4).Value = '2023-05-29'
$ws.Cells.Item(703, 5).Value = 5
$ws.Cells.Item(703, 6).Value = 100112037
$ws.Cells.Item(703, 7).Value = 'Cebollín'
$ws.Cells.Item(703, 8).Value = 'Sin especificar'
$ws.Cells.Item(703, 9).Value = 'Primera'
$ws.Cells.Item(703, 10).Value = 280
$ws.Cells.Item(703, 11).Value = 3500
$ws.Cells.Item(703, 12).Value = 4000
$ws.Cells.Item(703, 13).Value = 3786
$ws.Cells.Item(703, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(703, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(703, 16).Value = 105
$ws.Cells.Item(703, 17).Value = 36
$ws.Cells.Item(703, 18).Value = 'Hortaliza'

# New row 704: Cebollín, Segunda, same market date 2023-05-29
$ws.Cells.Item(704, 1).Value = 3
$ws.Cells.Item(704, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(704, 3).Value = 'Coquimbo'
$ws.Cells.Item(704, 4).Value = '2023-05-29'
$ws.Cells.Item(704, 5).Value = 5
$ws.Cells.Item(704, 6).Value = 100112037
$ws.Cells.Item(704, 7).Value = 'Cebollín'
$ws.Cells.Item(704, 8).Value = 'Sin especificar'
$ws.Cells.Item(704, 9).Value = 'Segunda'
$ws.Cells.Item(704, 10).Value = 110
$ws.Cells.Item(704, 11).Value = 3000
$ws.Cells.Item(704, 12).Value = 3000
$ws.Cells.Item(704, 13).Value = 3000
$ws.Cells.Item(704, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(704, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(704, 16).Value = 83
$ws.Cells.Item(704, 17).Value = 36
$ws.Cells.Item(704, 18).Value = 'Hortaliza'
